$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.159.92"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "2.940.98"
$ws.Range("E3").Value = "  -2.64%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "2.937.72"
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "65.161.52"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "3.431.60"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "2.940.42"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.16%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.90%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.32%  "
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.56%  "
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "385.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "2.703.98"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E50").Value = "  +4.37%  "
